# Auto-applied data update from scheduled runner
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 8
$ws.Range("H8").Value = 1431.909
$ws.Range("I8").Value = 1431.909
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 4295.727000000001
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -4156.727000000001
$ws.Range("N8").ClearContents()

# ALC row 107
$ws.Range("H107").Value = 359.75
$ws.Range("I107").Value = 306.08694
$ws.Range("J107").Value = 606.6
$ws.Range("K107").Value = 306.08694
$ws.Range("L107").Value = 606.6
$ws.Range("M107").Value = 1613.91306
$ws.Range("N107").Value = -4446.6

# ALC row 137
$ws.Range("H137").Value = 43480.527
$ws.Range("I137").Value = 1963.8182
$ws.Range("J137").Value = 61747.88
$ws.Range("K137").Value = 5891.4546
$ws.Range("L137").Value = 185243.64
$ws.Range("M137").Value = -3341.4546
$ws.Range("N137").Value = -190343.64

# ALC row 138
$ws.Range("H138").Value = 4313.0835
$ws.Range("I138").Value = 3677
$ws.Range("J138").Value = 4440.3
$ws.Range("K138").Value = 11031
$ws.Range("L138").Value = 13320.9
$ws.Range("M138").Value = -5891
$ws.Range("N138").Value = -23600.9

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 136765.28
$ws.Range("I32").Value = 141535.28
$ws.Range("J32").Value = 113975.336
$ws.Range("K32").Value = 141535.28
$ws.Range("L32").Value = 113975.336
$ws.Range("M32").Value = -141248.28
$ws.Range("N32").Value = -114549.336

# ARM row 61
$ws.Range("H61").Value = 2256.2856
$ws.Range("I61").Value = 2160.7307
$ws.Range("J61").Value = 3498.5
$ws.Range("K61").Value = 2160.7307
$ws.Range("L61").Value = 3498.5
$ws.Range("M61").Value = -1948.7307

# ARM row 74
$ws.Range("H74").Value = 1776.7
$ws.Range("I74").Value = 1369.2667
$ws.Range("J74").Value = 2999
$ws.Range("K74").Value = 1369.2667
$ws.Range("L74").Value = 2999
$ws.Range("M74").Value = -495.2666999999999
$ws.Range("N74").Value = -4747

# ARM row 77
$ws.Range("H77").Value = 1776.7
$ws.Range("I77").Value = 1369.2667
$ws.Range("J77").Value = 2999
$ws.Range("K77").Value = 6846.3335
$ws.Range("L77").Value = 14995
$ws.Range("M77").Value = -2478.3335
$ws.Range("N77").Value = -23731

# ARM row 97
$ws.Range("H97").Value = 1496.0834
$ws.Range("I97").Value = 1120.45
$ws.Range("J97").Value = 3374.25
$ws.Range("K97").Value = 1120.45
$ws.Range("L97").Value = 3374.25
$ws.Range("M97").Value = -624.45

# ARM row 102
$ws.Range("H102").Value = 2383.4167
$ws.Range("I102").Value = 2383.4167
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2383.4167
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -761.4167000000002

# ARM row 135
$ws.Range("H135").Value = 82664.664
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 82664.664
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 82664.664
$ws.Range("N135").Value = -92804.664

# ARM row 136
$ws.Range("H136").Value = 2256.2856
$ws.Range("I136").Value = 2160.7307
$ws.Range("J136").Value = 3498.5
$ws.Range("K136").Value = 6482.1921
$ws.Range("L136").Value = 10495.5
$ws.Range("M136").Value = -3932.1921

# ARM row 139
$ws.Range("H139").Value = 187995
$ws.Range("I139").Value = 187995
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 187995
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -182855

$ws = $wb.Worksheets.Item("BSM")
# BSM row 86
$ws.Range("H86").Value = 1994.0667
$ws.Range("I86").Value = 1633.3334
$ws.Range("J86").Value = 2535.1667
$ws.Range("K86").Value = 1633.3334
$ws.Range("L86").Value = 2535.1667
$ws.Range("M86").Value = -510.3334
$ws.Range("N86").Value = -4781.1667

# BSM row 89
$ws.Range("H89").Value = 1994.0667
$ws.Range("I89").Value = 1633.3334
$ws.Range("J89").Value = 2535.1667
$ws.Range("K89").Value = 8166.666999999999
$ws.Range("L89").Value = 12675.8335
$ws.Range("M89").Value = -2550.666999999999
$ws.Range("N89").Value = -23907.8335

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Range("H31").Value = 3202.8333
$ws.Range("I31").Value = 2679.2083
$ws.Range("J31").Value = 3551.9167
$ws.Range("K31").Value = 2679.2083
$ws.Range("L31").Value = 3551.9167
$ws.Range("M31").Value = -2384.2083
$ws.Range("N31").Value = -4141.9167

# CRP row 34
$ws.Range("H34").Value = 3202.8333
$ws.Range("I34").Value = 2679.2083
$ws.Range("J34").Value = 3551.9167
$ws.Range("K34").Value = 2679.2083
$ws.Range("L34").Value = 3551.9167
$ws.Range("M34").Value = -2477.2083
$ws.Range("N34").Value = -3955.9167

$ws = $wb.Worksheets.Item("CUL")
# CUL row 68
$ws.Range("H68").Value = 4171584
$ws.Range("I68").Value = 4144
$ws.Range("J68").Value = 10006000
$ws.Range("K68").Value = 12432
$ws.Range("L68").Value = 30018000
$ws.Range("M68").Value = -11621
$ws.Range("N68").Value = -30019622

# CUL row 71
$ws.Range("H71").Value = 4171584
$ws.Range("I71").Value = 4144
$ws.Range("J71").Value = 10006000
$ws.Range("K71").Value = 37296
$ws.Range("L71").Value = 90054000
$ws.Range("M71").Value = -33240
$ws.Range("N71").Value = -90062112

# CUL row 112
$ws.Range("H112").Value = 101117
$ws.Range("I112").Value = 101117
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 303351
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = -302243
$ws.Range("N112").ClearContents()

# CUL row 122
$ws.Range("H122").Value = 281.21738
$ws.Range("I122").Value = 145
$ws.Range("J122").Value = 309.89474
$ws.Range("K122").Value = 1305
$ws.Range("L122").Value = 2789.05266
$ws.Range("M122").Value = 1145
$ws.Range("N122").Value = -7689.05266

$ws = $wb.Worksheets.Item("GSM")
# GSM row 121
$ws.Range("H121").Value = 197281.8
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 197281.8
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 197281.8
$ws.Range("N121").Value = -200775.8

$ws = $wb.Worksheets.Item("LTW")
# LTW row 69
$ws.Range("H69").Value = 60000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 60000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 60000
$ws.Range("N69").Value = -61622

# LTW row 72
$ws.Range("H72").Value = 60000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 60000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 180000
$ws.Range("N72").Value = -188112

$ws = $wb.Worksheets.Item("WVR")
# WVR row 115
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("M115").ClearContents()

# WVR row 126
$ws.Range("H126").Value = 2326.8572
$ws.Range("I126").Value = 2382.1667
$ws.Range("J126").Value = 1995
$ws.Range("K126").Value = 7146.500100000001
$ws.Range("L126").Value = 5985
$ws.Range("M126").Value = -4676.500100000001
$ws.Range("N126").Value = -10925

# WVR row 129
$ws.Range("H129").Value = 100429
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 100429
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 100429
$ws.Range("N129").Value = -110429

# WVR row 132
$ws.Range("H132").Value = 7609.1562
$ws.Range("I132").Value = 8791.416999999999
$ws.Range("J132").Value = 4062.375
$ws.Range("K132").Value = 26374.251
$ws.Range("L132").Value = 12187.125
$ws.Range("M132").Value = -23844.251

# WVR row 136
$ws.Range("H136").Value = 6957
$ws.Range("I136").Value = 7134.909
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 21404.727
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -18854.727
